# Update Cmas/Etws xlsx file
#
# Net content changes on Sheet1 (values keep using the shared-string table;
# the engine manages dedup/garbage-collection of that table automatically):
#   - B2: "Tsunami warning" -> "Earthquake and tsunami warning"
#   - A3: 3 -> 2
#   - B3: "Earthquake and tsunami warning" -> "Tsunami warning"
#   - C3: long ETWS msg3.1 text -> "this is a ETWS test message2,this is a ETWS test message2.1"
#   - A4: 4 -> 3
#   - B4: "Earthquake warning" -> "Earthquake and tsunami warning"
#     (this removes the last reference to the "Earthquake warning" string)
#   - Selection moves from C20 to C18

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Earthquake and tsunami warning"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Tsunami warning"
$ws.Range("C3").Value = "this is a ETWS test message2,this is a ETWS test message2.1"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Earthquake and tsunami warning"

$ws.Range("C18").Select()
